$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$ws.Range('D2').Value = '23.161.06'
$ws.Range('E2').Value = '  -3.11%  '
$ws.Range('D3').Value = '1.607.24'
$ws.Range('E3').Value = '  -2.58%  '
$ws.Range('E4').Value = '  -0.07%  '
Set-TextValue $ws.Range('D5') '1.000'
$ws.Range('E5').Value = '  -0.03%  '
Set-TextValue $ws.Range('D6') '301.79'
$ws.Range('E6').Value = '  -2.18%  '
$ws.Range('E7').Value = '  -2.68%  '
Set-TextValue $ws.Range('D8') '0.3656'
$ws.Range('E8').Value = '  -4.46%  '
Set-TextValue $ws.Range('D9') '49.30'
$ws.Range('E9').Value = '  -5.37%  '
Set-TextValue $ws.Range('D10') '1.001'
$ws.Range('E10').Value = '  -0.01%  '
Set-TextValue $ws.Range('D11') '1.267'
$ws.Range('E11').Value = '  -6.30%  '
Set-TextValue $ws.Range('D12') '0.08089'
$ws.Range('E12').Value = '  -3.96%  '
Set-TextValue $ws.Range('D13') '23.05'
$ws.Range('E13').Value = '  -3.38%  '
Set-TextValue $ws.Range('D14') '6.610'
$ws.Range('E14').Value = '  -6.65%  '
Set-TextValue $ws.Range('D15') '7.447'
$ws.Range('E15').Value = '  -6.38%  '
$ws.Range('E16').Value = '  -4.63%  '
$ws.Range('D17').Value = '1.611.56'
$ws.Range('E17').Value = '  -2.21%  '
Set-TextValue $ws.Range('D18') '91.59'
$ws.Range('E18').Value = '  -3.26%  '
Set-TextValue $ws.Range('D19') '0.06802'
$ws.Range('E19').Value = '  -2.36%  '
Set-TextValue $ws.Range('D20') '18.39'
$ws.Range('E20').Value = '  -6.57%  '
Set-TextValue $ws.Range('D21') '6.579'
$ws.Range('E21').Value = '  -5.16%  '
Set-TextValue $ws.Range('D22') '1.001'
$ws.Range('E22').Value = '  +0.06%  '
Set-TextValue $ws.Range('D23') '13.09'
$ws.Range('E23').Value = '  -4.63%  '
$ws.Range('D24').Value = '23.186.22'
$ws.Range('E25').Value = '  -4.15%  '
Set-TextValue $ws.Range('D26') '2.874'
$ws.Range('E26').Value = '  -2.93%  '
$ws.Range('E27').Value = '  -4.46%  '
Set-TextValue $ws.Range('D28') '150.54'
$ws.Range('E28').Value = '  -0.68%  '
Set-TextValue $ws.Range('D29') '5.293'
$ws.Range('E29').Value = '  -2.00%  '
Set-TextValue $ws.Range('D30') '133.05'
$ws.Range('E30').Value = '  -4.32%  '
Set-TextValue $ws.Range('D31') '2.421'
$ws.Range('E31').Value = '  -3.87%  '
Set-TextValue $ws.Range('D32') '6.868'
$ws.Range('E32').Value = '  -12.82%  '
$ws.Range('D33').Value = '1.790.21'
$ws.Range('E33').Value = '  -2.09%  '
Set-TextValue $ws.Range('D34') '0.9712'
$ws.Range('E34').Value = '  -6.82%  '
Set-TextValue $ws.Range('D35') '0.07708'
$ws.Range('E35').Value = '  -4.38%  '
Set-TextValue $ws.Range('D36') '0.02767'
$ws.Range('E36').Value = '  -6.69%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D37') '6.265'
$ws.Range('E37').Value = '  -5.88%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D38') '0.2564'
$ws.Range('E38').Value = '  -4.38%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D39') '0.08906'
$ws.Range('E39').Value = '  -2.07%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D40') '10.12'
$ws.Range('E40').Value = '  -7.59%  '
Set-TextValue $ws.Range('D41') '1.391'
$ws.Range('E41').Value = '  -2.64%  '
Set-TextValue $ws.Range('D42') '0.7208'
$ws.Range('E42').Value = '  -5.42%  '
Set-TextValue $ws.Range('D43') '12.82'
$ws.Range('E43').Value = '  -4.67%  '
Set-TextValue $ws.Range('D44') '15.66'
$ws.Range('E44').Value = '  -4.51%  '
Set-TextValue $ws.Range('D45') '0.6688'
$ws.Range('E45').Value = '  -4.54%  '
Set-TextValue $ws.Range('D46') '2.313'
$ws.Range('E46').Value = '  -6.24%  '
Set-TextValue $ws.Range('D47') '0.9991'
$ws.Range('E47').Value = '  -0.13%  '
Set-TextValue $ws.Range('D48') '3.981'
$ws.Range('E48').Value = '  -2.29%  '
Set-TextValue $ws.Range('D49') '0.08020'
$ws.Range('E49').Value = '  -3.46%  '
Set-TextValue $ws.Range('D50') '130.78'
$ws.Range('E50').Value = '  -2.79%  '
$ws.Range('E51').Value = '  -2.89%  '
